$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.553.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.35%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.846.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.40%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.030"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +3.06%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'320.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.68%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.027"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.70%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4368"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.61%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3737"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.41%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07386"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +2.85%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.8742"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.80%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'21.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.51%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.856.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.89%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.483"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.11%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.673"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.37%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.07146"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.09%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'82.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.28%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.033"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.76%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000009008"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.48%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.027"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.62%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'15.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.36%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'27.565.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.35%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.240"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.35%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'11.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.92%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.066.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.04%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'157.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.33%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.926"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +4.60%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'18.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.51%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'5.260"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.34%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.950"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.47%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'115.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.65%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.09066"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.41%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.206"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.16%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.7647"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.04%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.503"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.45%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.872"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +4.19%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.028"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.55%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.144"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.34%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01976"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +3.36%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05261"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.77%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.5157"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.57%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.799"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +6.24%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1670"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.26%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'6.669"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +3.03%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.543"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.86%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'108.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.42%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'10.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.91%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.709"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +3.79%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.4644"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.39%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.06368"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.18%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.881"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +4.78%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'39.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +6.36%  "
$ws.Range("E51").Style = "Normal"
